# Auto-generated edit script: updates Sheets per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3635.0454
$ws.Range("I43").Value = 4724.1816
$ws.Range("J43").Value = 2545.9092
$ws.Range("K43").Value = 4724.1816
$ws.Range("L43").Value = 2545.9092
$ws.Range("M43").Value = -4655.1816
$ws.Range("N43").Value = -2683.9092

$ws.Range("H53").Value = 47620244
$ws.Range("I53").Value = 100001224
$ws.Range("K53").Value = 100001224
$ws.Range("M53").Value = -100000587

$ws.Range("H103").Value = 27778986
$ws.Range("I103").Value = 790
$ws.Range("J103").Value = 38462908
$ws.Range("K103").Value = 2370
$ws.Range("L103").Value = 115388724
$ws.Range("M103").Value = -1784
$ws.Range("N103").Value = -115389896

$ws.Range("H112").Value = 3637.3225
$ws.Range("J112").Value = 3637.3225
$ws.Range("L112").Value = 10911.9675
$ws.Range("N112").Value = -13127.9675

$ws.Range("H131").Value = 1749982.5
$ws.Range("I131").Value = 1875.0834
$ws.Range("J131").Value = 22727272
$ws.Range("K131").Value = 5625.2502
$ws.Range("L131").Value = 68181816
$ws.Range("M131").Value = -585.2502000000004
$ws.Range("N131").Value = -68191896

$ws.Range("H132").Value = 2503.3794
$ws.Range("I132").Value = 2629.7083
$ws.Range("K132").Value = 7889.124899999999
$ws.Range("M132").Value = -5359.124899999999

$ws.Range("H135").Value = 2778.8235
$ws.Range("J135").Value = 10500
$ws.Range("L135").Value = 94500
$ws.Range("N135").Value = -99570

$ws.Range("H137").Value = 2593.359
$ws.Range("I137").Value = 2420.6128
$ws.Range("J137").Value = 3262.75
$ws.Range("K137").Value = 7261.8384
$ws.Range("L137").Value = 9788.25
$ws.Range("M137").Value = -4711.8384
$ws.Range("N137").Value = -14888.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2799.5
$ws.Range("I2").Value = 2799.5
$ws.Range("K2").Value = 2799.5
$ws.Range("M2").Value = -2686.5

$ws.Range("H32").Value = 7506.7163
$ws.Range("I32").Value = 5166.5713
$ws.Range("J32").Value = 20909.363
$ws.Range("K32").Value = 5166.5713
$ws.Range("L32").Value = 20909.363
$ws.Range("M32").Value = -4879.5713
$ws.Range("N32").Value = -21483.363

$ws.Range("H45").Value = 3972139
$ws.Range("I45").Value = 15874891
$ws.Range("J45").Value = 4555.1113
$ws.Range("K45").Value = 15874891
$ws.Range("L45").Value = 4555.1113
$ws.Range("M45").Value = -15874514
$ws.Range("N45").Value = -5309.1113

$ws.Range("H61").Value = 3128931.8
$ws.Range("I61").Value = 3575246.5
$ws.Range("J61").Value = 4728
$ws.Range("K61").Value = 3575246.5
$ws.Range("L61").Value = 4728
$ws.Range("M61").Value = -3575034.5
$ws.Range("N61").Value = -5152

$ws.Range("H74").Value = 3754.7
$ws.Range("I74").Value = 3649.6667
$ws.Range("J74").Value = 3912.25
$ws.Range("K74").Value = 3649.6667
$ws.Range("L74").Value = 3912.25
$ws.Range("M74").Value = -2775.6667
$ws.Range("N74").Value = -5660.25

$ws.Range("H77").Value = 3754.7
$ws.Range("I77").Value = 3649.6667
$ws.Range("J77").Value = 3912.25
$ws.Range("K77").Value = 18248.3335
$ws.Range("L77").Value = 19561.25
$ws.Range("M77").Value = -13880.3335
$ws.Range("N77").Value = -28297.25

$ws.Range("H116").Value = 2799.5
$ws.Range("I116").Value = 2799.5
$ws.Range("K116").Value = 2799.5
$ws.Range("M116").Value = -505.5

$ws.Range("H132").Value = 4400.1763
$ws.Range("I132").Value = 4487.6875
$ws.Range("K132").Value = 13463.0625
$ws.Range("M132").Value = -10933.0625

$ws.Range("H136").Value = 3128931.8
$ws.Range("I136").Value = 3575246.5
$ws.Range("J136").Value = 4728
$ws.Range("K136").Value = 10725739.5
$ws.Range("L136").Value = 14184
$ws.Range("M136").Value = -10723189.5
$ws.Range("N136").Value = -19284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2799.5
$ws.Range("I3").Value = 2799.5
$ws.Range("K3").Value = 2799.5
$ws.Range("M3").Value = -2685.5

$ws.Range("H20").Value = 4552.5835
$ws.Range("I20").Value = 7080.8237
$ws.Range("J20").Value = 2290.4736
$ws.Range("K20").Value = 7080.8237
$ws.Range("L20").Value = 2290.4736
$ws.Range("M20").Value = -6833.8237
$ws.Range("N20").Value = -2784.4736

$ws.Range("H105").Value = 675763.7
$ws.Range("I105").Value = 1145429.6
$ws.Range("J105").Value = 4812.357
$ws.Range("K105").Value = 1145429.6
$ws.Range("L105").Value = 4812.357
$ws.Range("M105").Value = -1143682.6
$ws.Range("N105").Value = -8306.357

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24394678
$ws.Range("I31").Value = 71432390
$ws.Range("J31").Value = 4752.037
$ws.Range("K31").Value = 71432390
$ws.Range("L31").Value = 4752.037
$ws.Range("M31").Value = -71432095
$ws.Range("N31").Value = -5342.037

$ws.Range("H34").Value = 24394678
$ws.Range("I34").Value = 71432390
$ws.Range("J34").Value = 4752.037
$ws.Range("K34").Value = 71432390
$ws.Range("L34").Value = 4752.037
$ws.Range("M34").Value = -71432188
$ws.Range("N34").Value = -5156.037

$ws.Range("H58").Value = 2014.8206
$ws.Range("I58").Value = 1744.4667
$ws.Range("K58").Value = 1744.4667
$ws.Range("M58").Value = -1541.4667

$ws.Range("H86").Value = 4949.875
$ws.Range("I86").Value = 4317
$ws.Range("J86").Value = 5763.5713
$ws.Range("K86").Value = 4317
$ws.Range("L86").Value = 5763.5713
$ws.Range("M86").Value = -3194
$ws.Range("N86").Value = -8009.5713

$ws.Range("H89").Value = 4949.875
$ws.Range("I89").Value = 4317
$ws.Range("J89").Value = 5763.5713
$ws.Range("K89").Value = 21585
$ws.Range("L89").Value = 28817.8565
$ws.Range("M89").Value = -15969
$ws.Range("N89").Value = -40049.85649999999

$ws.Range("H102").Value = 80059.75
$ws.Range("I102").Value = 42000
$ws.Range("K102").Value = 42000
$ws.Range("M102").Value = -39566

$ws.Range("H103").Value = 43356
$ws.Range("J103").Value = 84999
$ws.Range("L103").Value = 84999
$ws.Range("N103").Value = -87343

$ws.Range("H104").Value = 25900
$ws.Range("I104").Value = 25900
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 25900
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = -23279
$ws.Range("N104").ClearContents()

$ws.Range("H107").Value = 1015.08
$ws.Range("I107").Value = 487.25
$ws.Range("J107").Value = 1953.4445
$ws.Range("K107").Value = 487.25
$ws.Range("L107").Value = 1953.4445
$ws.Range("M107").Value = 1432.75
$ws.Range("N107").Value = -5793.4445

$ws.Range("H132").Value = 1430.72
$ws.Range("I132").Value = 1399.0769
$ws.Range("J132").Value = 1542.909
$ws.Range("K132").Value = 4197.2307
$ws.Range("L132").Value = 4628.727000000001
$ws.Range("M132").Value = -1667.2307
$ws.Range("N132").Value = -9688.727000000001

$ws.Range("H136").Value = 2014.8206
$ws.Range("I136").Value = 1744.4667
$ws.Range("K136").Value = 5233.4001
$ws.Range("M136").Value = -2683.4001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 14407.538
$ws.Range("I56").Value = 14407.538
$ws.Range("K56").Value = 14407.538
$ws.Range("M56").Value = -13877.538

$ws.Range("H132").Value = 2834.1667
$ws.Range("I132").Value = 2416.1667
$ws.Range("J132").Value = 3252.1667
$ws.Range("K132").Value = 21745.5003
$ws.Range("L132").Value = 29269.5003
$ws.Range("M132").Value = -19215.5003
$ws.Range("N132").Value = -34329.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1356.28
$ws.Range("I102").Value = 1356.28
$ws.Range("K102").Value = 1356.28
$ws.Range("M102").Value = 265.72

$ws.Range("H109").Value = 99999.5
$ws.Range("J109").Value = 99999.5
$ws.Range("L109").Value = 99999.5
$ws.Range("N109").Value = -102079.5

$ws.Range("H132").Value = 2333554.8
$ws.Range("I132").Value = 2619.6316
$ws.Range("K132").Value = 7858.8948
$ws.Range("M132").Value = -5328.8948

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3369.8928
$ws.Range("I40").Value = 3204.875
$ws.Range("J40").Value = 4360
$ws.Range("K40").Value = 3204.875
$ws.Range("L40").Value = 4360
$ws.Range("M40").Value = -3068.875
$ws.Range("N40").Value = -4632

$ws.Range("H93").Value = 3707932.5
$ws.Range("I93").Value = 1690.3334
$ws.Range("K93").Value = 1690.3334
$ws.Range("M93").Value = -442.3334

$ws.Range("H102").Value = 99999
$ws.Range("J102").Value = 99999
$ws.Range("L102").Value = 99999
$ws.Range("N102").Value = -106489

$ws.Range("H132").Value = 2456.5066
$ws.Range("I132").Value = 1541.7142
$ws.Range("K132").Value = 4625.142599999999
$ws.Range("M132").Value = -2095.142599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2119.75
$ws.Range("I122").Value = 1960.5294
$ws.Range("J122").Value = 2506.4285
$ws.Range("K122").Value = 5881.5882
$ws.Range("L122").Value = 7519.2855
$ws.Range("M122").Value = -3431.5882
$ws.Range("N122").Value = -12419.2855

$ws.Range("H132").Value = 1278.3385
$ws.Range("I132").Value = 1148.2693
$ws.Range("J132").Value = 1798.6154
$ws.Range("K132").Value = 3444.8079
$ws.Range("L132").Value = 5395.8462
$ws.Range("M132").Value = -914.8078999999998
$ws.Range("N132").Value = -10455.8462

$ws.Range("H136").Value = 2489.1482
$ws.Range("I136").Value = 2384.6047
$ws.Range("J136").Value = 2897.818
$ws.Range("K136").Value = 7153.8141
$ws.Range("L136").Value = 8693.454000000002
$ws.Range("M136").Value = -4603.8141
$ws.Range("N136").Value = -13793.454

